$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-18 from 2023-09-02 to 2023-09-03
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45172
}
